# Commit: "Thu, May 28, 2020  7:05:56 PM"
#
# The deck's Design/Theme was switched from the custom "Integral" theme to
# PowerPoint's built-in default "Office Theme": the 10 non-neutral colours
# of the colour scheme that backs the (single) slide master / theme part
# change from the Integral palette to the stock Office palette (dk1=black
# and lt1=white are unchanged).
#
# We apply this with the documented PowerPoint COM surface: every Slide
# (and the SlideMaster) exposes a ThemeColorScheme of 12 colour slots
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) tied to the theme part
# backing the deck's one-and-only slide master. Setting .RGB on each slot
# rewrites that shared theme's <a:clrScheme> in place - exactly the edit
# recorded in the diff - without touching any slide content.

$p = $ppt.ActivePresentation

$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
